# Update NATMI LR-pair data (Cxcl10-Sdc4) with new TPM-derived values.
# Ligand-side (G,H,I,J) values are keyed by the Sending cluster (column A);
# Receptor-side (M,N,O,P) values are keyed by the Target cluster (column D);
# Edge weights (Q,R,S,T) are simply the products G*M, H*N, I*O, J*P.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ligand (sending-cluster) derived values: avg expr, total expr, specificity(avg), specificity(total)
$ligand = @{
    "ECs"               = @(31.35623066666667,  94.068692,  0.07215642027787079, 0.07299614919666826)
    "FAPs"              = @(139.9983773333333, 419.995132,  0.3221618650682612,  0.3259110620709639)
    "Inflammatory-Mac"  = @(130.001713,         390.005139,  0.2991577125385502,  0.3026391959814998)
    "MuSCs"             = @(14.997169,          29.994338,   0.03451122811430998, 0.02327523775607825)
    "Resolving-Mac"     = @(118.2056323333333,  354.616897,  0.2720127740010079,  0.2751783549947898)
}

# New receptor (target-cluster) derived values: avg expr, total expr, specificity(avg), specificity(total)
$receptor = @{
    "ECs"               = @(1.378421333333333, 4.135264,            0.01656231489052403, 0.01794267551419991)
    "FAPs"              = @(18.067884,          54.20365200000001,  0.2170932623988173,  0.2351865659654651)
    "Inflammatory-Mac"  = @(17.58286933333333,  52.748608,           0.2112656061941426,  0.22887321273073)
    "MuSCs"             = @(19.2082395,         38.416479,           0.2307951156866419,  0.1666869194070983)
    "Resolving-Mac"     = @(26.988955,          80.966865,           0.3242837008298742,  0.3513106263825066)
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $sender = $ws.Cells.Item($r, 1).Value2   # column A - Sending cluster
    $target = $ws.Cells.Item($r, 4).Value2   # column D - Target cluster

    if (-not $ligand.ContainsKey($sender)) { continue }
    if (-not $receptor.ContainsKey($target)) { continue }

    $ligVals = $ligand[$sender]
    $recVals = $receptor[$target]

    $g = $ligVals[0]
    $h = $ligVals[1]
    $i = $ligVals[2]
    $j = $ligVals[3]

    $m = $recVals[0]
    $n = $recVals[1]
    $o = $recVals[2]
    $p = $recVals[3]

    $ws.Cells.Item($r, 7).Value  = $g   # G
    $ws.Cells.Item($r, 8).Value  = $h   # H
    $ws.Cells.Item($r, 9).Value  = $i   # I
    $ws.Cells.Item($r, 10).Value = $j   # J

    $ws.Cells.Item($r, 13).Value = $m   # M
    $ws.Cells.Item($r, 14).Value = $n   # N
    $ws.Cells.Item($r, 15).Value = $o   # O
    $ws.Cells.Item($r, 16).Value = $p   # P

    $ws.Cells.Item($r, 17).Value = $g * $m   # Q
    $ws.Cells.Item($r, 18).Value = $h * $n   # R
    $ws.Cells.Item($r, 19).Value = $i * $o   # S
    $ws.Cells.Item($r, 20).Value = $j * $p   # T
}
